# completed E5 titrations 0331
# Adds the new 2022-03-31 CRM accuracy titration row (row 75) to the
# CRMAccuracyData sheet, matching the existing table layout/formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

$newRow = 75

$ws.Cells.Item($newRow, 1).Value = 20220331
$ws.Cells.Item($newRow, 2).Value = 2226.165
$ws.Cells.Item($newRow, 3).Value = 2224.4699999999998
$ws.Cells.Item($newRow, 4).Formula = "=100*(B" + $newRow + "-C" + $newRow + ")/C" + $newRow
$ws.Cells.Item($newRow, 5).Value = 180
$ws.Cells.Item($newRow, 6).Value = "CRM OPENED 20220318"

# Move/update the sheet's active selection to mirror where the user
# continued entry after adding this row.
$ws.Range("B76").Select()
